$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 3).Value = 2.834502210481444
$ws.Cells.Item(2, 4).Value = 2.921953217054471
$ws.Cells.Item(2, 5).Value = 16.7027395269962
$ws.Cells.Item(2, 6).Value = 20.05797417919059
$ws.Cells.Item(2, 7).Value = 24.47440444073665
$ws.Cells.Item(2, 8).Value = 11.41363832979199
$ws.Cells.Item(2, 9).Value = 16.38611221003454
$ws.Cells.Item(2, 14).Value = 18.99769123883711
$ws.Cells.Item(2, 15).Value = 17.07954586972049

# Row 3
$ws.Cells.Item(3, 3).Value = 2.808503672603281
$ws.Cells.Item(3, 4).Value = 2.921686010819876
$ws.Cells.Item(3, 5).Value = 15.73979365400252
$ws.Cells.Item(3, 6).Value = 19.55191849635492
$ws.Cells.Item(3, 7).Value = 23.37815233946195
$ws.Cells.Item(3, 8).Value = 11.34212789423975
$ws.Cells.Item(3, 9).Value = 16.09835956784373
$ws.Cells.Item(3, 14).Value = 18.40031508502701
$ws.Cells.Item(3, 15).Value = 16.76986830385361

# Row 4
$ws.Cells.Item(4, 3).Value = 2.7931115702881
$ws.Cells.Item(4, 4).Value = 2.921759718579023
$ws.Cells.Item(4, 5).Value = 15.12273035622041
$ws.Cells.Item(4, 6).Value = 19.24103833372577
$ws.Cells.Item(4, 7).Value = 22.68586246961432
$ws.Cells.Item(4, 8).Value = 11.30089585243082
$ws.Cells.Item(4, 9).Value = 15.9241860918018
$ws.Cells.Item(4, 14).Value = 18.02485520896362
$ws.Cells.Item(4, 15).Value = 16.58213466568255

# Row 5
$ws.Cells.Item(5, 3).Value = 2.786989497203994
$ws.Cells.Item(5, 4).Value = 2.921849622287124
$ws.Cells.Item(5, 5).Value = 14.86504626354035
$ws.Cells.Item(5, 6).Value = 19.11452490033642
$ws.Cells.Item(5, 7).Value = 22.39939137988523
$ws.Cells.Item(5, 8).Value = 11.28478012596251
$ws.Cells.Item(5, 9).Value = 15.85394126840663
$ws.Cells.Item(5, 14).Value = 17.86990355188765
$ws.Cells.Item(5, 15).Value = 16.50634811581754

# Row 6
$ws.Cells.Item(6, 3).Value = 2.785982202230878
$ws.Cells.Item(6, 4).Value = 2.921868168734245
$ws.Cells.Item(6, 5).Value = 14.82189077943935
$ws.Cells.Item(6, 6).Value = 19.09353423208228
$ws.Cells.Item(6, 7).Value = 22.35157473898746
$ws.Cells.Item(6, 8).Value = 11.28214597632855
$ws.Cells.Item(6, 9).Value = 15.84232435492578
$ws.Cells.Item(6, 14).Value = 17.84406337566597
$ws.Cells.Item(6, 15).Value = 16.49381030925665

# Row 7
$ws.Cells.Item(7, 3).Value = 2.793028388484726
$ws.Cells.Item(7, 4).Value = 2.92176068852762
$ws.Cells.Item(7, 5).Value = 15.11927995754622
$ws.Cells.Item(7, 6).Value = 19.23933114346418
$ws.Cells.Item(7, 7).Value = 22.68201600144635
$ws.Cells.Item(7, 8).Value = 11.30067571247396
$ws.Cells.Item(7, 9).Value = 15.92323564679213
$ws.Cells.Item(7, 14).Value = 18.02277304767603
$ws.Cells.Item(7, 15).Value = 16.58110953531916

# Row 8
$ws.Cells.Item(8, 3).Value = 2.825423502010025
$ws.Cells.Item(8, 4).Value = 2.921811781287785
$ws.Cells.Item(8, 5).Value = 16.37621570653944
$ws.Cells.Item(8, 6).Value = 19.88364487339634
$ws.Cells.Item(8, 7).Value = 24.10066096499387
$ws.Cells.Item(8, 8).Value = 11.38843486064322
$ws.Cells.Item(8, 9).Value = 16.28643572825336
$ws.Cells.Item(8, 14).Value = 18.79364780656867
$ws.Cells.Item(8, 15).Value = 16.97233525539156

# Row 9
$ws.Cells.Item(9, 3).Value = 2.893173323068783
$ws.Cells.Item(9, 4).Value = 2.923795484154702
$ws.Cells.Item(9, 5).Value = 18.7794094180861
$ws.Cells.Item(9, 6).Value = 21.13692313316237
$ws.Cells.Item(9, 7).Value = 26.71277809985466
$ws.Cells.Item(9, 8).Value = 11.58109830804097
$ws.Cells.Item(9, 9).Value = 17.01409898298131
$ws.Cells.Item(9, 14).Value = 20.2273683202997
$ws.Cells.Item(9, 15).Value = 17.75380677738767

# Row 10
$ws.Cells.Item(10, 3).Value = 2.945095738912221
$ws.Cells.Item(10, 4).Value = 2.926397385555783
$ws.Cells.Item(10, 5).Value = 20.4563019867387
$ws.Cells.Item(10, 6).Value = 22.04018808408795
$ws.Cells.Item(10, 7).Value = 28.50850979450653
$ws.Cells.Item(10, 8).Value = 11.73422960676105
$ws.Cells.Item(10, 9).Value = 17.5522062056261
$ws.Cells.Item(10, 14).Value = 21.22223697909767
$ws.Cells.Item(10, 15).Value = 18.3302793736581

# Row 11
$ws.Cells.Item(11, 3).Value = 2.969088735992326
$ws.Cells.Item(11, 4).Value = 2.927828707778827
$ws.Cells.Item(11, 5).Value = 21.17649355917162
$ws.Cells.Item(11, 6).Value = 22.44514234115031
$ws.Cells.Item(11, 7).Value = 29.29539143617675
$ws.Cells.Item(11, 8).Value = 11.80618519361414
$ws.Cells.Item(11, 9).Value = 17.79655404344218
$ws.Cells.Item(11, 14).Value = 21.66018057919901
$ws.Cells.Item(11, 15).Value = 18.59173846750071

# Row 12
$ws.Cells.Item(12, 3).Value = 2.978219953220289
$ws.Cells.Item(12, 4).Value = 2.928406236134329
$ws.Cells.Item(12, 5).Value = 21.44311882955938
$ws.Cells.Item(12, 6).Value = 22.5974700366945
$ws.Cells.Item(12, 7).Value = 29.5888309626042
$ws.Cells.Item(12, 8).Value = 11.83374313777817
$ws.Cells.Item(12, 9).Value = 17.88892260391525
$ws.Cells.Item(12, 14).Value = 21.82377585682186
$ws.Cells.Item(12, 15).Value = 18.69053090742778

# Row 13
$ws.Cells.Item(13, 3).Value = 2.976251483190431
$ws.Cells.Item(13, 4).Value = 2.928280277095107
$ws.Cells.Item(13, 5).Value = 21.38596685349934
$ws.Cells.Item(13, 6).Value = 22.56471137841228
$ws.Cells.Item(13, 7).Value = 29.52583837487112
$ws.Cells.Item(13, 8).Value = 11.8277945843078
$ws.Cells.Item(13, 9).Value = 17.86903803392805
$ws.Cells.Item(13, 14).Value = 21.78864458690801
$ws.Cells.Item(13, 15).Value = 18.66926542142354

# Row 14
$ws.Cells.Item(14, 3).Value = 2.969839099305567
$ws.Cells.Item(14, 4).Value = 2.927875510157183
$ws.Cells.Item(14, 5).Value = 21.198550945294
$ws.Cells.Item(14, 6).Value = 22.45769576608209
$ws.Cells.Item(14, 7).Value = 29.31962492154876
$ws.Cells.Item(14, 8).Value = 11.80844630836669
$ws.Cells.Item(14, 9).Value = 17.80415702362767
$ws.Cells.Item(14, 14).Value = 21.67368539489659
$ws.Cells.Item(14, 15).Value = 18.59987109836093

# Row 15
$ws.Cells.Item(15, 3).Value = 2.965917026405468
$ws.Cells.Item(15, 4).Value = 2.927632201064195
$ws.Cells.Item(15, 5).Value = 21.08296049945123
$ws.Cells.Item(15, 6).Value = 22.39200815117349
$ws.Cells.Item(15, 7).Value = 29.19271661739445
$ws.Cells.Item(15, 8).Value = 11.79663468415387
$ws.Cells.Item(15, 9).Value = 17.76439178787451
$ws.Cells.Item(15, 14).Value = 21.60297336126124
$ws.Cells.Item(15, 15).Value = 18.55733387493612

# Row 16
$ws.Cells.Item(16, 3).Value = 2.943534559880475
$ws.Cells.Item(16, 4).Value = 2.926308821137225
$ws.Cells.Item(16, 5).Value = 20.40838086130212
$ws.Cells.Item(16, 6).Value = 22.01358935488452
$ws.Cells.Item(16, 7).Value = 28.45646176678661
$ws.Cells.Item(16, 8).Value = 11.72957162029912
$ws.Cells.Item(16, 9).Value = 17.53622022891942
$ws.Cells.Item(16, 14).Value = 21.19330956972086
$ws.Cells.Item(16, 15).Value = 18.31316768169064

# Row 17
$ws.Cells.Item(17, 3).Value = 2.929893489230503
$ws.Cells.Item(17, 4).Value = 2.925560342612977
$ws.Cells.Item(17, 5).Value = 19.98365333694488
$ws.Cells.Item(17, 6).Value = 21.77979758343491
$ws.Cells.Item(17, 7).Value = 27.99693971637666
$ws.Cells.Item(17, 8).Value = 11.68900436561275
$ws.Cells.Item(17, 9).Value = 17.39606113090706
$ws.Cells.Item(17, 14).Value = 20.93814219015166
$ws.Cells.Item(17, 15).Value = 18.16310401839639

# Row 18
$ws.Cells.Item(18, 3).Value = 2.922083121824703
$ws.Cells.Item(18, 4).Value = 2.925153159089977
$ws.Cells.Item(18, 5).Value = 19.73534616995746
$ws.Cells.Item(18, 6).Value = 21.64477615587819
$ws.Cells.Item(18, 7).Value = 27.72982129512856
$ws.Cells.Item(18, 8).Value = 11.66588845484039
$ws.Cells.Item(18, 9).Value = 17.315409512441
$ws.Cells.Item(18, 14).Value = 20.79000725568362
$ws.Cells.Item(18, 15).Value = 18.07672379090914

# Row 19
$ws.Cells.Item(19, 3).Value = 2.919445028707286
$ws.Cells.Item(19, 4).Value = 2.925019302266418
$ws.Cells.Item(19, 5).Value = 19.65058295046888
$ws.Cells.Item(19, 6).Value = 21.59897088845631
$ws.Cells.Item(19, 7).Value = 27.63890367331787
$ws.Cells.Item(19, 8).Value = 11.6580997302975
$ws.Cells.Item(19, 9).Value = 17.28809914378477
$ws.Cells.Item(19, 14).Value = 20.73962067985785
$ws.Cells.Item(19, 15).Value = 18.04746857153799

# Row 20
$ws.Cells.Item(20, 3).Value = 2.931341976901293
$ws.Cells.Item(20, 4).Value = 2.925637606425657
$ws.Cells.Item(20, 5).Value = 20.02928160063097
$ws.Cells.Item(20, 6).Value = 21.80474333456365
$ws.Cells.Item(20, 7).Value = 28.04614951869
$ws.Cells.Item(20, 8).Value = 11.69330047348328
$ws.Cells.Item(20, 9).Value = 17.41098573139981
$ws.Cells.Item(20, 14).Value = 20.96544799484619
$ws.Cells.Item(20, 15).Value = 18.17908631470507

# Row 21
$ws.Cells.Item(21, 3).Value = 2.971721398969928
$ws.Cells.Item(21, 4).Value = 2.927993436917687
$ws.Cells.Item(21, 5).Value = 21.25376469477846
$ws.Cells.Item(21, 6).Value = 22.48915776837041
$ws.Cells.Item(21, 7).Value = 29.38031949041861
$ws.Cells.Item(21, 8).Value = 11.81412111816867
$ws.Cells.Item(21, 9).Value = 17.82321926515937
$ws.Cells.Item(21, 14).Value = 21.70751365554066
$ws.Cells.Item(21, 15).Value = 18.62026057364821

# Row 22
$ws.Cells.Item(22, 3).Value = 2.998373683119186
$ws.Cells.Item(22, 4).Value = 2.929740081760873
$ws.Cells.Item(22, 5).Value = 22.01852129484435
$ws.Cells.Item(22, 6).Value = 22.93045133276988
$ws.Cells.Item(22, 7).Value = 30.22576161766649
$ws.Cells.Item(22, 8).Value = 11.89488114835638
$ws.Cells.Item(22, 9).Value = 18.09165908564879
$ws.Cells.Item(22, 14).Value = 22.17935961385674
$ws.Cells.Item(22, 15).Value = 18.90728699567935

# Row 23
$ws.Cells.Item(23, 3).Value = 2.984127550908692
$ws.Cells.Item(23, 4).Value = 2.928788961440894
$ws.Cells.Item(23, 5).Value = 21.61359423408213
$ws.Cells.Item(23, 6).Value = 22.69552572357221
$ws.Cells.Item(23, 7).Value = 29.77702354447638
$ws.Cells.Item(23, 8).Value = 11.85162039647262
$ws.Cells.Item(23, 9).Value = 17.94850833355131
$ws.Cells.Item(23, 14).Value = 21.92877110911181
$ws.Cells.Item(23, 15).Value = 18.75424828378117

# Row 24
$ws.Cells.Item(24, 3).Value = 2.930687015147974
$ws.Cells.Item(24, 4).Value = 2.92560260343284
$ws.Cells.Item(24, 5).Value = 20.00866591621054
$ws.Cells.Item(24, 6).Value = 21.79346725221375
$ws.Cells.Item(24, 7).Value = 28.02391088581589
$ws.Cells.Item(24, 8).Value = 11.69135755756665
$ws.Cells.Item(24, 9).Value = 17.40423853716608
$ws.Cells.Item(24, 14).Value = 20.95310750188673
$ws.Cells.Item(24, 15).Value = 18.17186104263799

# Row 25
$ws.Cells.Item(25, 3).Value = 2.874436613818586
$ws.Cells.Item(25, 4).Value = 2.923057563585686
$ws.Cells.Item(25, 5).Value = 18.12383623505977
$ws.Cells.Item(25, 6).Value = 20.80018057611101
$ws.Cells.Item(25, 7).Value = 26.0265544071983
$ws.Cells.Item(25, 8).Value = 11.5268708537711
$ws.Cells.Item(25, 9).Value = 16.81621279524714
$ws.Cells.Item(25, 14).Value = 19.84905939529497
$ws.Cells.Item(25, 15).Value = 17.5415371771969

Write-Output "Updated loading_percent values for rows 2-25 (380 kV case)"